$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 250
$ws.Range("I2").Value = 242.22223
$ws.Range("J2").Value = 285
$ws.Range("K2").Value = 242.22223
$ws.Range("L2").Value = 285
$ws.Range("M2").Value = -129.22223
$ws.Range("N2").Value = -511

$ws.Range("H21").Value = 50000
$ws.Range("J21").Value = 50000
$ws.Range("L21").Value = 50000
$ws.Range("N21").Value = -50936

$ws.Range("H23").Value = 50000
$ws.Range("J23").Value = 50000
$ws.Range("L23").Value = 50000
$ws.Range("N23").Value = -50468

$ws.Range("H38").Value = 224.33333
$ws.Range("I38").Value = 62.909092
$ws.Range("J38").Value = 2000
$ws.Range("K38").Value = 188.727276
$ws.Range("L38").Value = 6000
$ws.Range("M38").Value = 183.272724
$ws.Range("N38").Value = -6744

$ws.Range("H58").Value = 299
$ws.Range("I58").Value = 111.25
$ws.Range("K58").Value = 333.75
$ws.Range("M58").Value = -183.75

$ws.Range("H87").Value = 16688.271
$ws.Range("J87").Value = 16688.271
$ws.Range("L87").Value = 16688.271
$ws.Range("N87").Value = -19184.271

$ws.Range("H90").Value = 16688.271
$ws.Range("J90").Value = 16688.271
$ws.Range("L90").Value = 50064.813
$ws.Range("N90").Value = -62544.813

$ws.Range("H116").Value = 1823.5
$ws.Range("I116").Value = 1892.7778
$ws.Range("J116").Value = 1200
$ws.Range("K116").Value = 1892.7778
$ws.Range("L116").Value = 1200
$ws.Range("M116").Value = 1549.2222
$ws.Range("N116").Value = -8084

$ws.Range("H129").Value = 844.2234999999999
$ws.Range("J129").Value = 894.24
$ws.Range("L129").Value = 2682.72
$ws.Range("N129").Value = -12682.72

$ws.Range("H141").Value = 5139.7
$ws.Range("I141").Value = 1488.5555
$ws.Range("J141").Value = 38000
$ws.Range("K141").Value = 4465.666499999999
$ws.Range("L141").Value = 114000
$ws.Range("M141").Value = 714.3335000000006
$ws.Range("N141").Value = -124360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 55701.43
$ws.Range("J140").Value = 55701.43
$ws.Range("L140").Value = 55701.43
$ws.Range("N140").Value = -66061.42999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H138").Value = 39758.8
$ws.Range("J138").Value = 39758.8
$ws.Range("L138").Value = 39758.8
$ws.Range("N138").Value = -50038.8

$ws.Range("H140").Value = 80054.5
$ws.Range("J140").Value = 80054.5
$ws.Range("L140").Value = 80054.5
$ws.Range("N140").Value = -90414.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1347.375
$ws.Range("I97").Value = 1340
$ws.Range("J97").Value = 1349.8334
$ws.Range("K97").Value = 4020
$ws.Range("L97").Value = 4049.5002
$ws.Range("M97").Value = -3524
$ws.Range("N97").Value = -5041.5002

$ws.Range("H118").Value = 867.5789
$ws.Range("I118").Value = 387.1111
$ws.Range("J118").Value = 1300
$ws.Range("K118").Value = 1161.3333
$ws.Range("L118").Value = 3900
$ws.Range("M118").Value = 81.66669999999999
$ws.Range("N118").Value = -6386

$ws.Range("H131").Value = 907.46
$ws.Range("J131").Value = 923.19794
$ws.Range("L131").Value = 2769.59382
$ws.Range("N131").Value = -12849.59382

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1858.909
$ws.Range("I113").Value = 1695
$ws.Range("J113").Value = 2296
$ws.Range("K113").Value = 1695
$ws.Range("L113").Value = 2296
$ws.Range("M113").Value = 475
$ws.Range("N113").Value = -6636

$ws.Range("H122").Value = 106517.125
$ws.Range("I122").Value = 179898.28
$ws.Range("J122").Value = 3783.5
$ws.Range("K122").Value = 539694.84
$ws.Range("L122").Value = 11350.5
$ws.Range("M122").Value = -537244.84
$ws.Range("N122").Value = -16250.5

$ws.Range("H132").Value = 2789.5833
$ws.Range("I132").Value = 2534.7693
$ws.Range("J132").Value = 3090.7273
$ws.Range("K132").Value = 7604.3079
$ws.Range("L132").Value = 9272.1819
$ws.Range("M132").Value = -5074.3079
$ws.Range("N132").Value = -14332.1819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 103690.3
$ws.Range("I7").Value = 169870.5
$ws.Range("J7").Value = 4420
$ws.Range("K7").Value = 169870.5
$ws.Range("L7").Value = 4420
$ws.Range("M7").Value = -169758.5
$ws.Range("N7").Value = -4644

$ws.Range("H46").Value = 1688.95
$ws.Range("I46").Value = 1248.25
$ws.Range("M46").Value = -1060.25

$ws.Range("H68").Value = 1613.5834
$ws.Range("I68").Value = 1045
$ws.Range("K68").Value = 1045
$ws.Range("M68").Value = -296
$ws.Range("N68").Value = -3225.3

$ws.Range("H71").Value = 1613.5834
$ws.Range("I71").Value = 1045
$ws.Range("K71").Value = 5225
$ws.Range("L71").Value = 8636.5
$ws.Range("M71").Value = -1481
$ws.Range("N71").Value = -16124.5

$ws.Range("H82").Value = 1969.3529
$ws.Range("I82").Value = 1855.6428
$ws.Range("J82").Value = 2500
$ws.Range("K82").Value = 1855.6428
$ws.Range("L82").Value = 2500
$ws.Range("M82").Value = -1494.6428
$ws.Range("N82").Value = -3222

$ws.Range("H85").Value = 1969.3529
$ws.Range("I85").Value = 1855.6428
$ws.Range("J85").Value = 2500
$ws.Range("K85").Value = 1855.6428
$ws.Range("L85").Value = 2500
$ws.Range("M85").Value = -607.6428000000001
$ws.Range("N85").Value = -4996

$ws.Range("H126").Value = 103690.3
$ws.Range("I126").Value = 169870.5
$ws.Range("J126").Value = 4420
$ws.Range("K126").Value = 509611.5
$ws.Range("L126").Value = 13260
$ws.Range("M126").Value = -507141.5
$ws.Range("N126").Value = -18200

$ws.Range("H127").Value = 49548.125
$ws.Range("J127").Value = 49548.125
$ws.Range("L127").Value = 49548.125
$ws.Range("N127").Value = -59468.125

$ws.Range("H133").Value = 47646
$ws.Range("J133").Value = 47646
$ws.Range("L133").Value = 47646
$ws.Range("N133").Value = -52706

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3100
$ws.Range("I62").Value = 3350
$ws.Range("J62").Value = 2800
$ws.Range("K62").Value = 3350
$ws.Range("L62").Value = 2800
$ws.Range("M62").Value = -2726
$ws.Range("N62").Value = -4048

$ws.Range("H65").Value = 3100
$ws.Range("I65").Value = 3350
$ws.Range("J65").Value = 2800
$ws.Range("K65").Value = 16750
$ws.Range("L65").Value = 14000
$ws.Range("M65").Value = -13630
$ws.Range("N65").Value = -20240

$ws.Range("H136").Value = 3188
$ws.Range("I136").Value = 629.8
$ws.Range("J136").Value = 5663.6772
$ws.Range("K136").Value = 1889.4
$ws.Range("L136").Value = 16991.0316
$ws.Range("M136").Value = 660.6000000000001
$ws.Range("N136").Value = -22091.0316
